$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '89.769.32'
$ws.Range("E2").Value = '  -0.93%  '

$ws.Range("D3").Value = '3.070.58'
$ws.Range("E3").Value = '  -2.19%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.74'
$ws.Range("E5").Value = '  +10.01%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '616.88'
$ws.Range("E6").Value = '  -0.94%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.06'
$ws.Range("E7").Value = '  -6.39%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.366'
$ws.Range("E8").Value = '  +0.58%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  +0.05%  '

$ws.Range("D10").Value = '3.070.54'
$ws.Range("E10").Value = '  -2.14%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.705'
$ws.Range("E11").Value = '  -4.67%  '

$ws.Range("E12").Value = '  -0.43%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000248'
$ws.Range("E13").Value = '  +0.99%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.67'
$ws.Range("E14").Value = '  -1.63%  '

$ws.Range("D15").Value = '89.349.82'
$ws.Range("E15").Value = '  -1.03%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.36'
$ws.Range("E16").Value = '  -5.24%  '

$ws.Range("D17").Value = '3.637.48'
$ws.Range("E17").Value = '  -2.01%  '

$ws.Range("D18").Value = '3.058.77'
$ws.Range("E18").Value = '  -2.55%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.78'
$ws.Range("E19").Value = '  +0.56%  '

$ws.Range("E20").Value = '  -0.05%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.75'
$ws.Range("E21").Value = '  -5.67%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '430.91'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.41'
$ws.Range("E23").Value = '  +1.00%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.69'
$ws.Range("E24").Value = '  -4.29%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.70'
$ws.Range("E25").Value = '  -1.17%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '86.93'
$ws.Range("E26").Value = '  -8.52%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.68'
$ws.Range("E27").Value = '  -5.09%  '

$ws.Range("E28").Value = '  -2.32%  '

$ws.Range("E29").Value = '  +0.20%  '

$ws.Range("E30").Value = '  +41.35%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.158'
$ws.Range("E31").Value = '  -3.50%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.94'
$ws.Range("E32").Value = '  -2.98%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.198'
$ws.Range("E33").Value = '  -7.82%  '

$ws.Range("B34").Value = 'MantraDAO'
$ws.Range("C34").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.16'
$ws.Range("E34").Value = '  +65.08%  '

$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '25.56'
$ws.Range("E35").Value = '  -4.43%  '

$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.150'
$ws.Range("E36").Value = '  +2.76%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.12'
$ws.Range("E37").Value = '  +1.47%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '491.05'
$ws.Range("E38").Value = '  -5.31%  '

$ws.Range("E39").Value = '  +0.20%  '

$ws.Range("E40").Value = '  -3.22%  '

$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.25'
$ws.Range("E41").Value = '  -5.85%  '

$ws.Range("B42").Value = 'Hedera'
$ws.Range("C42").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0899'
$ws.Range("E42").Value = '  -1.70%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.08'
$ws.Range("E43").Value = '  -0.62%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.398'
$ws.Range("E45").Value = '  -6.85%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '156.48'
$ws.Range("E46").Value = '  +4.13%  '

$ws.Range("E47").Value = '  -6.92%  '

$ws.Range("E48").Value = '  -8.01%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '44.46'
$ws.Range("E49").Value = '  -1.99%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("E50").Value = '  -0.10%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.30'
$ws.Range("E51").Value = '  -5.45%  '
